$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 4
$ws.Range("C1").Formula = "=SUM(A1,B1)"
$ws.Range("D1").Formula = "=SUM(A1,B1)/7+A1"

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 5

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 6

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4

$ws.Range("B5").Select()
